$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A128").Value = "Login with valid username and password"
$ws.Range("B128").Value = "PASSED"
$ws.Range("C128").Value = "chrome"
$ws.Range("A129").Value = "Create a country"
$ws.Range("B129").Value = "PASSED"
$ws.Range("C129").Value = "chrome"
$ws.Range("A130").Value = "Create a country"
$ws.Range("B130").Value = "PASSED"
$ws.Range("C130").Value = "chrome"
$ws.Range("A131").Value = "Create a country"
$ws.Range("B131").Value = "PASSED"
$ws.Range("C131").Value = "chrome"
$ws.Range("A132").Value = "Create a country 2"
$ws.Range("B132").Value = "PASSED"
$ws.Range("C132").Value = "chrome"
$ws.Range("A133").Value = "Login with valid username and password"
$ws.Range("B133").Value = "PASSED"
$ws.Range("C133").Value = "chrome"
$ws.Range("A134").Value = "Create a country"
$ws.Range("B134").Value = "PASSED"
$ws.Range("C134").Value = "chrome"
$ws.Range("A135").Value = "Create a country"
$ws.Range("B135").Value = "PASSED"
$ws.Range("C135").Value = "chrome"
$ws.Range("A136").Value = "Create a country 2"
$ws.Range("B136").Value = "PASSED"
$ws.Range("C136").Value = "chrome"
$ws.Range("A137").Value = "Create a citizenship"
$ws.Range("B137").Value = "PASSED"
$ws.Range("C137").Value = "chrome"
$ws.Range("A138").Value = "Create a Citizenship"
$ws.Range("B138").Value = "PASSED"
$ws.Range("C138").Value = "chrome"
$ws.Range("A139").Value = "Create a Citizenship"
$ws.Range("B139").Value = "PASSED"
$ws.Range("C139").Value = "chrome"
$ws.Range("A140").Value = "Create a Citizenship"
$ws.Range("B140").Value = "PASSED"
$ws.Range("C140").Value = "chrome"
$ws.Range("A141").Value = "Create a Citizenship"
$ws.Range("B141").Value = "PASSED"
$ws.Range("C141").Value = "chrome"
$ws.Range("A142").Value = "Create a Citizenship"
$ws.Range("B142").Value = "PASSED"
$ws.Range("C142").Value = "chrome"
$ws.Range("A143").Value = "Create Country"
$ws.Range("B143").Value = "FAILED"
$ws.Range("C143").Value = "chrome"
$ws.Range("A144").Value = "Create Nationality"
$ws.Range("B144").Value = "PASSED"
$ws.Range("C144").Value = "chrome"
$ws.Range("A145").Value = "Fee Functionality"
$ws.Range("B145").Value = "PASSED"
$ws.Range("C145").Value = "chrome"
$ws.Range("A146").Value = "Fee Functionality"
$ws.Range("B146").Value = "PASSED"
$ws.Range("C146").Value = "chrome"
$ws.Range("A147").Value = "Fee Functionality"
$ws.Range("B147").Value = "PASSED"
$ws.Range("C147").Value = "chrome"
$ws.Range("A148").Value = "Fee Functionality"
$ws.Range("B148").Value = "PASSED"
$ws.Range("C148").Value = "chrome"
$ws.Range("A149").Value = "Fee Functionality"
$ws.Range("B149").Value = "PASSED"
$ws.Range("C149").Value = "chrome"
$ws.Range("A150").Value = "Fee Functionality"
$ws.Range("B150").Value = "PASSED"
$ws.Range("C150").Value = "chrome"
$ws.Range("A151").Value = "Create Nationality and Delete"
$ws.Range("B151").Value = "PASSED"
$ws.Range("C151").Value = "chrome"
$ws.Range("A152").Value = "Create Nationality and Delete"
$ws.Range("B152").Value = "PASSED"
$ws.Range("C152").Value = "chrome"
$ws.Range("A153").Value = "Create Nationality and Delete"
$ws.Range("B153").Value = "PASSED"
$ws.Range("C153").Value = "chrome"
$ws.Range("A154").Value = "Create Nationality and Delete"
$ws.Range("B154").Value = "PASSED"
$ws.Range("C154").Value = "chrome"
$ws.Range("A155").Value = "Create Nationality and Delete"
$ws.Range("B155").Value = "PASSED"
$ws.Range("C155").Value = "chrome"
$ws.Range("A156").Value = "Create Inventory and Delete"
$ws.Range("B156").Value = "PASSED"
$ws.Range("C156").Value = "chrome"
$ws.Range("A157").Value = "Create Inventory and Delete"
$ws.Range("B157").Value = "PASSED"
$ws.Range("C157").Value = "chrome"
$ws.Range("A158").Value = "Create Inventory and Delete"
$ws.Range("B158").Value = "PASSED"
$ws.Range("C158").Value = "chrome"
$ws.Range("A159").Value = "Create Inventory and Delete"
$ws.Range("B159").Value = "PASSED"
$ws.Range("C159").Value = "chrome"
$ws.Range("A160").Value = "Create Inventory and Delete"
$ws.Range("B160").Value = "PASSED"
$ws.Range("C160").Value = "chrome"
$ws.Range("A161").Value = "Create and Delete CitizenShip From Excel"
$ws.Range("B161").Value = "PASSED"
$ws.Range("C161").Value = "chrome"
$ws.Range("A162").Value = "States testing with JDBC"
$ws.Range("B162").Value = "FAILED"
$ws.Range("C162").Value = "chrome"
